# Update the five "problem rows" of the divide-and-check table with the
# newly generated two-digit/one-digit division problems.
#
# The table layout is unchanged (20 rows x 5 columns, with data only in
# rows 1, 5, 9, 13 and 17 and the remaining rows left blank for students
# to work in) -- only the cell text in those five rows is replaced.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$rows = @{
    1  = @("57÷9=6, 3",  "56÷7=8, 0",  "71÷7=10, 1", "60÷4=15, 0", "18÷6=3, 0")
    5  = @("28÷7=4, 0",  "18÷2=9, 0",  "74÷8=9, 2",  "27÷7=3, 6",  "69÷5=13, 4")
    9  = @("26÷8=3, 2",  "37÷4=9, 1",  "83÷2=41, 1", "59÷9=6, 5",  "86÷6=14, 2")
    13 = @("53÷7=7, 4",  "35÷3=11, 2", "49÷7=7, 0",  "25÷4=6, 1",  "11÷7=1, 4")
    17 = @("99÷3=33, 0", "71÷6=11, 5", "84÷6=14, 0", "32÷8=4, 0",  "77÷2=38, 1")
}

foreach ($rowIndex in $rows.Keys) {
    $values = $rows[$rowIndex]
    for ($col = 1; $col -le $values.Length; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
